# Auto-generated edit script applying numeric corrections to Leve profit tables
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1196.7778
$ws.Range("I4").Value = 1196.7778
$ws.Range("K4").Value = 1196.7778
$ws.Range("M4").Value = -1082.7778
$ws.Range("H86").Value = 83370950
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 90950030
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 90950030
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -90952276
$ws.Range("H89").Value = 83370950
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 90950030
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 454750150
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -454761382
$ws.Range("H98").Value = 32889.39
$ws.Range("J98").Value = 16040.75
$ws.Range("L98").Value = 16040.75
$ws.Range("N98").Value = -19036.75
$ws.Range("H106").Value = 10444.538
$ws.Range("I106").Value = 10098.125
$ws.Range("K106").Value = 10098.125
$ws.Range("M106").Value = -9467.125
$ws.Range("H122").Value = 32889.39
$ws.Range("J122").Value = 16040.75
$ws.Range("L122").Value = 48122.25
$ws.Range("N122").Value = -53022.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 78588.41
$ws.Range("I45").Value = 147400
$ws.Range("K45").Value = 147400
$ws.Range("M45").Value = -147023
$ws.Range("H61").Value = 7972.3667
$ws.Range("I61").Value = 9282.362999999999
$ws.Range("J61").Value = 4369.875
$ws.Range("K61").Value = 9282.362999999999
$ws.Range("L61").Value = 4369.875
$ws.Range("M61").Value = -9070.362999999999
$ws.Range("N61").Value = -4793.875
$ws.Range("H136").Value = 7972.3667
$ws.Range("I136").Value = 9282.362999999999
$ws.Range("J136").Value = 4369.875
$ws.Range("K136").Value = 27847.089
$ws.Range("L136").Value = 13109.625
$ws.Range("M136").Value = -25297.089
$ws.Range("N136").Value = -18209.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 797
$ws.Range("I5").Value = 299
$ws.Range("K5").Value = 299
$ws.Range("M5").Value = -186
$ws.Range("H20").Value = 3607.889
$ws.Range("I20").Value = 2674.625
$ws.Range("J20").Value = 4354.5
$ws.Range("K20").Value = 2674.625
$ws.Range("L20").Value = 4354.5
$ws.Range("M20").Value = -2427.625
$ws.Range("N20").Value = -4848.5
$ws.Range("H99").Value = 10908.4375
$ws.Range("I99").Value = 11890.6
$ws.Range("K99").Value = 11890.6
$ws.Range("M99").Value = -10392.6
$ws.Range("H134").Value = 15877.333
$ws.Range("I134").Value = 17584.25
$ws.Range("K134").Value = 52752.75
$ws.Range("M134").Value = -50217.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1640.238
$ws.Range("I22").Value = 1190.5
$ws.Range("J22").Value = 1820.1333
$ws.Range("K22").Value = 1190.5
$ws.Range("L22").Value = 1820.1333
$ws.Range("M22").Value = -840.5
$ws.Range("N22").Value = -2520.1333
$ws.Range("H58").Value = 3748
$ws.Range("I58").Value = 3469.3076
$ws.Range("J58").Value = 4472.6
$ws.Range("K58").Value = 3469.3076
$ws.Range("L58").Value = 4472.6
$ws.Range("M58").Value = -3266.3076
$ws.Range("N58").Value = -4878.6
$ws.Range("H134").Value = 2536.5715
$ws.Range("I134").Value = 2522.8696
$ws.Range("K134").Value = 7568.6088
$ws.Range("M134").Value = -5033.6088
$ws.Range("H135").Value = 119555
$ws.Range("J135").Value = 119555
$ws.Range("L135").Value = 119555
$ws.Range("N135").Value = -129695
$ws.Range("H136").Value = 3748
$ws.Range("I136").Value = 3469.3076
$ws.Range("J136").Value = 4472.6
$ws.Range("K136").Value = 10407.9228
$ws.Range("L136").Value = 13417.8
$ws.Range("M136").Value = -7857.9228
$ws.Range("N136").Value = -18517.8
$ws.Range("H138").Value = 62997.5
$ws.Range("J138").Value = 62997.5
$ws.Range("L138").Value = 62997.5
$ws.Range("N138").Value = -73277.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 173.33333
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 173.33333
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 519.99999
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -799.99999
$ws.Range("H14").Value = 423.5
$ws.Range("I14").Value = 423.5
$ws.Range("K14").Value = 1270.5
$ws.Range("M14").Value = -1097.5
$ws.Range("H46").Value = 4519.1665
$ws.Range("I46").Value = 1710.4
$ws.Range("J46").Value = 5599.4614
$ws.Range("K46").Value = 5131.200000000001
$ws.Range("L46").Value = 16798.3842
$ws.Range("M46").Value = -5040.200000000001
$ws.Range("N46").Value = -16980.3842
$ws.Range("H128").Value = 339999
$ws.Range("I128").Value = 339999
$ws.Range("K128").Value = 1019997
$ws.Range("M128").Value = -1015017

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 1787.5
$ws.Range("I19").Value = 1775
$ws.Range("J19").Value = 1800
$ws.Range("K19").Value = 1775
$ws.Range("L19").Value = 1800
$ws.Range("M19").Value = -1487
$ws.Range("N19").Value = -2376
$ws.Range("H80").Value = 11088.857
$ws.Range("I80").Value = 10670.454
$ws.Range("K80").Value = 10670.454
$ws.Range("M80").Value = -9672.454
$ws.Range("H83").Value = 11088.857
$ws.Range("I83").Value = 10670.454
$ws.Range("K83").Value = 53352.27
$ws.Range("M83").Value = -48360.27
$ws.Range("H126").Value = 9403.5
$ws.Range("I126").Value = 14211.8
$ws.Range("K126").Value = 42635.39999999999
$ws.Range("M126").Value = -40165.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6703.5454
$ws.Range("I16").Value = 6703.5454
$ws.Range("K16").Value = 6703.5454
$ws.Range("M16").Value = -6533.5454
$ws.Range("H20").Value = 12995
$ws.Range("I20").Value = 12995
$ws.Range("K20").Value = 12995
$ws.Range("M20").Value = -12769
$ws.Range("H22").Value = 5970.55
$ws.Range("I22").Value = 12949.625
$ws.Range("J22").Value = 1317.8334
$ws.Range("K22").Value = 12949.625
$ws.Range("L22").Value = 1317.8334
$ws.Range("M22").Value = -12654.625
$ws.Range("N22").Value = -1907.8334
$ws.Range("H27").Value = 5970.55
$ws.Range("I27").Value = 12949.625
$ws.Range("J27").Value = 1317.8334
$ws.Range("K27").Value = 12949.625
$ws.Range("L27").Value = 1317.8334
$ws.Range("M27").Value = -12842.625
$ws.Range("N27").Value = -1531.8334
$ws.Range("H93").Value = 10969.1
$ws.Range("I93").Value = 11211.5
$ws.Range("K93").Value = 11211.5
$ws.Range("M93").Value = -9963.5
$ws.Range("H122").Value = 6596.3335
$ws.Range("I122").Value = 6596.3335
$ws.Range("K122").Value = 19789.0005
$ws.Range("M122").Value = -17339.0005
$ws.Range("H135").Value = 94760.78999999999
$ws.Range("J135").Value = 94760.78999999999
$ws.Range("L135").Value = 94760.78999999999
$ws.Range("N135").Value = -104900.79

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4443.24
$ws.Range("I122").Value = 1906.5
$ws.Range("J122").Value = 7191.375
$ws.Range("K122").Value = 5719.5
$ws.Range("L122").Value = 21574.125
$ws.Range("M122").Value = -3269.5
$ws.Range("N122").Value = -26474.125
$ws.Range("H132").Value = 7822
$ws.Range("I132").Value = 9580.106
$ws.Range("K132").Value = 28740.318
$ws.Range("M132").Value = -26210.318
